# Rename the two logo pictures that are inlined into the document's
# headers and footers. The BTec logo (JPG), embedded in both the
# primary and first-page headers, moves from "image1.jpg" to
# "image2.jpg". The Pearson/Edexcel logo (PNG), embedded in both the
# primary and first-page footers, moves from "image2.png" to
# "image1.png".

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Headers: index 1 = primary ("default"), index 2 = first-page header.
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}

# Footers: index 1 = primary ("default"), index 2 = first-page footer.
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -like "*PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}
